$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("B9").Value = 90857

# Row 10
$ws.Range("A10").Value = 112344211
$ws.Range("B10").Value = 90837
$ws.Range("E10").Value = 5966
$ws.Range("F10").Value = "Motaggsvamp"
$ws.Range("G10").Value = "Sarcodon squamosus"
$ws.Range("H10").Value = "(Schaeff.) Quél."
$ws.Range("Q10").Value = 763527
$ws.Range("R10").Value = 7089456
$ws.Range("Z10").Value = "12:26"
$ws.Range("AB10").Value = "12:26"

# Row 11
$ws.Range("B11").Value = 90814

# Row 12
$ws.Range("A12").Value = 112344251
$ws.Range("B12").Value = 90857
$ws.Range("E12").Value = 5448
$ws.Range("F12").Value = "Svartvit taggsvamp"
$ws.Range("G12").Value = "Phellodon connatus"
$ws.Range("H12").Value = "(Schultz) nom.prov"
$ws.Range("Q12").Value = 763530
$ws.Range("R12").Value = 7089425
$ws.Range("Z12").Value = "12:29"
$ws.Range("AB12").Value = "12:29"
